# Rename the SKU prefix from "PIPI" to "BEAT" for every data row in column A
# (what used to be a nested loop in the old robot/RPA flow is now a single
# pass over the sheet's used range in this COM script).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $sku = $cell.Value2
    if ($sku -like "PIPI*") {
        $cell.Value = "BEAT" + $sku.Substring(4)
    }
}
